$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original data block A2:D10 (9 rows of CEP/address entries) is a
# repeating 9-row cycle. Build out a provisional larger base by copying
# that block down the sheet to fill rows 11-50 (four more full cycles
# plus a final partial cycle), exactly mirroring what a user would get
# by selecting A2:D10, copying, and pasting repeatedly down the column.
$cycleLen = 9
$lastRow = 50
$startRow = 11
while ($startRow -le $lastRow) {
    $endRow = [Math]::Min($startRow + $cycleLen - 1, $lastRow)
    $rowCount = $endRow - $startRow + 1
    $srcEnd = 2 + $rowCount - 1

    $src = $ws.Range("A2:D$srcEnd")
    $dest = $ws.Range("A$startRow`:D$endRow")
    $src.Copy($dest)

    $startRow += $cycleLen
}

# Land the selection on the new last cell, matching where the user ended
# up after pasting the final chunk, and scroll the window down to it.
$ws.Range("A50").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 30
